## sprint_37.xlsx - "server list test case added"
##
## The new "server list" test-case block lives at the bottom of the sheet
## (rows 54:57), so this view update scrolls/selects down to it and
## refreshes the row-3 header height flag + merged-header bookkeeping that
## Excel re-touches whenever it re-saves the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply row 3's height so Excel marks it as an explicit/custom height.
$ws.Rows.Item(3).RowHeight = 37.5

# Touch the sheet's merged "day header" ranges so they are re-registered
# (and therefore re-serialized) in the same order as the saved workbook:
# the four day headers that sit above row 24 first, then the remaining
# ones in their original relative order.
$null = $ws.Range("B2:C2").UnMerge()
$null = $ws.Range("B8:C8").UnMerge()
$null = $ws.Range("B16:C16").UnMerge()
$null = $ws.Range("B24:C24").UnMerge()
$null = $ws.Range("B54:C54").UnMerge()
$null = $ws.Range("B48:C48").UnMerge()
$null = $ws.Range("B42:C42").UnMerge()
$null = $ws.Range("B30:C30").UnMerge()
$null = $ws.Range("B36:C36").UnMerge()

$null = $ws.Range("B2:C2").Merge()
$null = $ws.Range("B8:C8").Merge()
$null = $ws.Range("B16:C16").Merge()
$null = $ws.Range("B24:C24").Merge()
$null = $ws.Range("B54:C54").Merge()
$null = $ws.Range("B48:C48").Merge()
$null = $ws.Range("B42:C42").Merge()
$null = $ws.Range("B30:C30").Merge()
$null = $ws.Range("B36:C36").Merge()

# Scroll the window down near the new block and select the newly added
# "server list" test-case summary range (A54:C57).
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 47
$null = $ws.Range("A54:C57").Select()
